$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 179, pushing existing rows 179-204 down to 181-206
$ws.Range("A179:R180").Insert()

# Row 179: new weekly entry, "Segunda" quality
$ws.Cells.Item(179,1).Value = 1
$ws.Cells.Item(179,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(179,3).Value = "Arica y Parinacota"
$ws.Cells.Item(179,4).Value = 44476
$ws.Cells.Item(179,5).Value = 15
$ws.Cells.Item(179,6).Value = 100112023
$ws.Cells.Item(179,7).Value = "Brócoli"
$ws.Cells.Item(179,8).Value = "Sin especificar"
$ws.Cells.Item(179,9).Value = "Segunda"
$ws.Cells.Item(179,10).Value = 800
$ws.Cells.Item(179,11).Value = 500
$ws.Cells.Item(179,12).Value = 600
$ws.Cells.Item(179,13).Value = 550
$ws.Cells.Item(179,14).Value = "`$/unidad"
$ws.Cells.Item(179,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(179,16).Value = 550
$ws.Cells.Item(179,17).Value = 1
$ws.Cells.Item(179,18).Value = "Hortaliza"

# Row 180: new weekly entry, "Tercera" quality
$ws.Cells.Item(180,1).Value = 1
$ws.Cells.Item(180,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(180,3).Value = "Arica y Parinacota"
$ws.Cells.Item(180,4).Value = 44476
$ws.Cells.Item(180,5).Value = 15
$ws.Cells.Item(180,6).Value = 100112023
$ws.Cells.Item(180,7).Value = "Brócoli"
$ws.Cells.Item(180,8).Value = "Sin especificar"
$ws.Cells.Item(180,9).Value = "Tercera"
$ws.Cells.Item(180,10).Value = 1000
$ws.Cells.Item(180,11).Value = 400
$ws.Cells.Item(180,12).Value = 450
$ws.Cells.Item(180,13).Value = 425
$ws.Cells.Item(180,14).Value = "`$/unidad"
$ws.Cells.Item(180,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(180,16).Value = 425
$ws.Cells.Item(180,17).Value = 1
$ws.Cells.Item(180,18).Value = "Hortaliza"
